# Updates cryptos list: prices and 1h volume % changes, plus a few
# coin-row position swaps, per the Dec 10 2023 scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.057.05"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.373.37"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.39"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.51"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.98"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +4.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.20"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +15.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.31"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.109"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.42"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.929"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.378.29"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.012.90"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.05"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.70"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.79"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.42"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.75"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.43"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.10"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.81"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  +6.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.41"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0283"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.56"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +18.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.72"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +10.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.89"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +14.31%  "
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.11"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.60"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.25"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.81"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  +17.71%  "

Write-Output "Updated cryptos list"
